# Update the 2014 "count" value (428 -> 432)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 432

# Add the new "2015" row: A6 = "2015" (stored as text, matching the other
# year cells in column A), B6 = 292
$yearCell = $ws.Range("A6")
$yearCell.NumberFormat = "@"
$yearCell.Value = "2015"
$yearCell.Style = "Normal"

$ws.Range("B6").Value = 292
